$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.069678564573214
$ws.Cells.Item(2, 4).Value = 1.070862449233158
$ws.Cells.Item(2, 5).Value = 1.073974098786153
$ws.Cells.Item(2, 6).Value = 1.083875362078025
$ws.Cells.Item(2, 9).Value = 1.055816520875321
$ws.Cells.Item(2, 10).Value = 1.074612070367332
$ws.Cells.Item(2, 11).Value = 1.073561155671785
$ws.Cells.Item(2, 12).Value = 1.076664539590814
$ws.Cells.Item(2, 13).Value = 1.086539842605904
$ws.Cells.Item(2, 14).Value = 1.076138142766981

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.070804053876516
$ws.Cells.Item(3, 4).Value = 1.071738165512357
$ws.Cells.Item(3, 5).Value = 1.074955216544574
$ws.Cells.Item(3, 6).Value = 1.084885415757727
$ws.Cells.Item(3, 9).Value = 1.056152564049975
$ws.Cells.Item(3, 10).Value = 1.07539397126052
$ws.Cells.Item(3, 11).Value = 1.074253168507703
$ws.Cells.Item(3, 12).Value = 1.077462290867877
$ws.Cells.Item(3, 13).Value = 1.087368335119355
$ws.Cells.Item(3, 14).Value = 1.076921154049123

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.071532604912567
$ws.Cells.Item(4, 4).Value = 1.072305035546442
$ws.Cells.Item(4, 5).Value = 1.075590598230983
$ws.Cells.Item(4, 6).Value = 1.085539566606131
$ws.Cells.Item(4, 9).Value = 1.056368984015565
$ws.Cells.Item(4, 10).Value = 1.075899594543636
$ws.Cells.Item(4, 11).Value = 1.074700527454672
$ws.Cells.Item(4, 12).Value = 1.077978399741001
$ws.Cells.Item(4, 13).Value = 1.087904385166167
$ws.Cells.Item(4, 14).Value = 1.077427495375297

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.071838956049593
$ws.Cells.Item(5, 4).Value = 1.072543400635415
$ws.Cells.Item(5, 5).Value = 1.075857839865926
$ws.Cells.Item(5, 6).Value = 1.085814709691299
$ws.Cells.Item(5, 9).Value = 1.05645972205821
$ws.Cells.Item(5, 10).Value = 1.076112082065859
$ws.Cells.Item(5, 11).Value = 1.074888496437916
$ws.Cells.Item(5, 12).Value = 1.07819535020311
$ws.Cells.Item(5, 13).Value = 1.088129730698998
$ws.Cells.Item(5, 14).Value = 1.077640284654173

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.071890397764809
$ws.Cells.Item(6, 4).Value = 1.072583426272726
$ws.Cells.Item(6, 5).Value = 1.075902718358299
$ws.Cells.Item(6, 6).Value = 1.08586091550955
$ws.Cells.Item(6, 9).Value = 1.056474943009227
$ws.Cells.Item(6, 10).Value = 1.076147755164744
$ws.Cells.Item(6, 11).Value = 1.074920051343732
$ws.Cells.Item(6, 12).Value = 1.078231775859559
$ws.Cells.Item(6, 13).Value = 1.088167566615103
$ws.Cells.Item(6, 14).Value = 1.077676008412951

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.071536698122831
$ws.Cells.Item(7, 4).Value = 1.072308220385046
$ws.Cells.Item(7, 5).Value = 1.075594168627578
$ws.Cells.Item(7, 6).Value = 1.085543242539535
$ws.Cells.Item(7, 9).Value = 1.056370197424228
$ws.Cells.Item(7, 10).Value = 1.075902434113124
$ws.Cells.Item(7, 11).Value = 1.074703039501495
$ws.Cells.Item(7, 12).Value = 1.077981298729325
$ws.Cells.Item(7, 13).Value = 1.087907396283538
$ws.Cells.Item(7, 14).Value = 1.077430338977299

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.070058870253229
$ws.Cells.Item(8, 4).Value = 1.071158355183948
$ws.Cells.Item(8, 5).Value = 1.074305561545868
$ws.Cells.Item(8, 6).Value = 1.084216594382454
$ws.Cells.Item(8, 9).Value = 1.055930299859339
$ws.Cells.Item(8, 10).Value = 1.074876383360762
$ws.Cells.Item(8, 11).Value = 1.073795111381874
$ws.Cells.Item(8, 12).Value = 1.076934161748267
$ws.Cells.Item(8, 13).Value = 1.086819843615922
$ws.Cells.Item(8, 14).Value = 1.076402831115179

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.067456914161877
$ws.Cells.Item(9, 4).Value = 1.069133874530885
$ws.Cells.Item(9, 5).Value = 1.07203896995672
$ws.Cells.Item(9, 6).Value = 1.081883321353081
$ws.Cells.Item(9, 9).Value = 1.055147319422088
$ws.Cells.Item(9, 10).Value = 1.073065915644923
$ws.Cells.Item(9, 11).Value = 1.072192023280503
$ws.Cells.Item(9, 12).Value = 1.075088298441252
$ws.Cells.Item(9, 13).Value = 1.084903144681253
$ws.Cells.Item(9, 14).Value = 1.07458979232752

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.065723714259715
$ws.Cells.Item(10, 4).Value = 1.067785408671057
$ws.Cells.Item(10, 5).Value = 1.070530685865825
$ws.Cells.Item(10, 6).Value = 1.080330822035578
$ws.Cells.Item(10, 9).Value = 1.054620078431268
$ws.Cells.Item(10, 10).Value = 1.071857307159609
$ws.Cells.Item(10, 11).Value = 1.071121157627599
$ws.Cells.Item(10, 12).Value = 1.073857283286744
$ws.Cells.Item(10, 13).Value = 1.083625166053566
$ws.Cells.Item(10, 14).Value = 1.073379467479527

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.064973555533089
$ws.Cells.Item(11, 4).Value = 1.067201794555617
$ws.Cells.Item(11, 5).Value = 1.06987824459546
$ws.Cells.Item(11, 6).Value = 1.07965929231854
$ws.Cells.Item(11, 9).Value = 1.054390532135079
$ws.Cells.Item(11, 10).Value = 1.07133358052988
$ws.Cells.Item(11, 11).Value = 1.070656955469638
$ws.Cells.Item(11, 12).Value = 1.073324138055048
$ws.Cells.Item(11, 13).Value = 1.083071747581215
$ws.Cells.Item(11, 14).Value = 1.072854997097911

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.064694961962275
$ws.Cells.Item(12, 4).Value = 1.066985056730675
$ws.Cells.Item(12, 5).Value = 1.069635997582188
$ws.Cells.Item(12, 6).Value = 1.079409963618138
$ws.Cells.Item(12, 9).Value = 1.054305081026735
$ws.Cells.Item(12, 10).Value = 1.071138986339907
$ws.Cells.Item(12, 11).Value = 1.070484453546394
$ws.Cells.Item(12, 12).Value = 1.073126088109804
$ws.Cells.Item(12, 13).Value = 1.082866176663262
$ws.Cells.Item(12, 14).Value = 1.072660126561869

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.064754719030055
$ws.Cells.Item(13, 4).Value = 1.067031545804385
$ws.Cells.Item(13, 5).Value = 1.069687955917921
$ws.Cells.Item(13, 6).Value = 1.079463440609411
$ws.Cells.Item(13, 9).Value = 1.054323419059068
$ws.Cells.Item(13, 10).Value = 1.07118073011554
$ws.Cells.Item(13, 11).Value = 1.070521459258161
$ws.Cells.Item(13, 12).Value = 1.073168571220999
$ws.Cells.Item(13, 13).Value = 1.082910272621801
$ws.Cells.Item(13, 14).Value = 1.072701929618451

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.064950525909587
$ws.Cells.Item(14, 4).Value = 1.067183878071827
$ws.Cells.Item(14, 5).Value = 1.06985821836658
$ws.Cells.Item(14, 6).Value = 1.079638680539267
$ws.Cells.Item(14, 9).Value = 1.054383472542422
$ws.Cells.Item(14, 10).Value = 1.071317496513115
$ws.Cells.Item(14, 11).Value = 1.070642697968401
$ws.Cells.Item(14, 12).Value = 1.073307767512045
$ws.Cells.Item(14, 13).Value = 1.083054755159931
$ws.Cells.Item(14, 14).Value = 1.072838890239998

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.065071175438501
$ws.Cells.Item(15, 4).Value = 1.06727774064253
$ws.Cells.Item(15, 5).Value = 1.069963135768957
$ws.Cells.Item(15, 6).Value = 1.079746665874519
$ws.Cells.Item(15, 9).Value = 1.054420448644877
$ws.Cells.Item(15, 10).Value = 1.071401755011999
$ws.Cells.Item(15, 11).Value = 1.070717386995158
$ws.Cells.Item(15, 12).Value = 1.073393528814593
$ws.Cells.Item(15, 13).Value = 1.08314377474265
$ws.Cells.Item(15, 14).Value = 1.072923268395615

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.065773506694038
$ws.Cells.Item(16, 4).Value = 1.067824147153291
$ws.Cells.Item(16, 5).Value = 1.070574000077108
$ws.Cells.Item(16, 6).Value = 1.080375404289559
$ws.Cells.Item(16, 9).Value = 1.054635286382269
$ws.Cells.Item(16, 10).Value = 1.071892056926794
$ws.Cells.Item(16, 11).Value = 1.07115195445281
$ws.Cells.Item(16, 12).Value = 1.073892664121626
$ws.Cells.Item(16, 13).Value = 1.083661893673275
$ws.Cells.Item(16, 14).Value = 1.073414266595368

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.06621414804114
$ws.Cells.Item(17, 4).Value = 1.068166968831764
$ws.Cells.Item(17, 5).Value = 1.070957354817301
$ws.Cells.Item(17, 6).Value = 1.080769986285526
$ws.Cells.Item(17, 9).Value = 1.054769714452949
$ws.Cells.Item(17, 10).Value = 1.072199505612133
$ws.Cells.Item(17, 11).Value = 1.07142441058746
$ws.Cells.Item(17, 12).Value = 1.074205729782315
$ws.Cells.Item(17, 13).Value = 1.083986883967893
$ws.Cells.Item(17, 14).Value = 1.073722151893105

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.066471198392762
$ws.Cells.Item(18, 4).Value = 1.068366957960727
$ws.Cells.Item(18, 5).Value = 1.071181022401966
$ws.Cells.Item(18, 6).Value = 1.081000208039038
$ws.Cells.Item(18, 9).Value = 1.054848003707872
$ws.Cells.Item(18, 10).Value = 1.072378797487364
$ws.Cells.Item(18, 11).Value = 1.071583280402138
$ws.Cells.Item(18, 12).Value = 1.074388325302753
$ws.Cells.Item(18, 13).Value = 1.084176440928144
$ws.Cells.Item(18, 14).Value = 1.073901698383363

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.066558851314697
$ws.Cells.Item(19, 4).Value = 1.068435153663507
$ws.Cells.Item(19, 5).Value = 1.071257297990153
$ws.Cells.Item(19, 6).Value = 1.081078719361117
$ws.Cells.Item(19, 9).Value = 1.054874677904034
$ws.Cells.Item(19, 10).Value = 1.072439924982884
$ws.Cells.Item(19, 11).Value = 1.071637442548928
$ws.Cells.Item(19, 12).Value = 1.074450583892729
$ws.Cells.Item(19, 13).Value = 1.084241074205555
$ws.Cells.Item(19, 14).Value = 1.073962912686939

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.066166868118728
$ws.Cells.Item(20, 4).Value = 1.068130184510554
$ws.Cells.Item(20, 5).Value = 1.070916217937774
$ws.Cells.Item(20, 6).Value = 1.080727644258458
$ws.Cells.Item(20, 9).Value = 1.054755304042808
$ws.Cells.Item(20, 10).Value = 1.072166523181752
$ws.Cells.Item(20, 11).Value = 1.071395183715026
$ws.Cells.Item(20, 12).Value = 1.074172141872956
$ws.Cells.Item(20, 13).Value = 1.083952016036866
$ws.Cells.Item(20, 14).Value = 1.07368912262389

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.064892864321616
$ws.Cells.Item(21, 4).Value = 1.067139018841791
$ws.Cells.Item(21, 5).Value = 1.069808077592124
$ws.Cells.Item(21, 6).Value = 1.079587073782587
$ws.Cells.Item(21, 9).Value = 1.054365793461131
$ws.Cells.Item(21, 10).Value = 1.071277223841259
$ws.Cells.Item(21, 11).Value = 1.070606998303105
$ws.Cells.Item(21, 12).Value = 1.073266778116529
$ws.Cells.Item(21, 13).Value = 1.083012208836128
$ws.Cells.Item(21, 14).Value = 1.072798560376329

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.064092129539551
$ws.Cells.Item(22, 4).Value = 1.066516079358336
$ws.Cells.Item(22, 5).Value = 1.069111917071581
$ws.Cells.Item(22, 6).Value = 1.078870572998807
$ws.Cells.Item(22, 9).Value = 1.054119808512717
$ws.Cells.Item(22, 10).Value = 1.070717745833849
$ws.Cells.Item(22, 11).Value = 1.070110991868304
$ws.Cells.Item(22, 12).Value = 1.072697447071348
$ws.Cells.Item(22, 13).Value = 1.082421276311718
$ws.Cells.Item(22, 14).Value = 1.072238287845976

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.064516587648516
$ws.Cells.Item(23, 4).Value = 1.06684628797351
$ws.Cells.Item(23, 5).Value = 1.069480910711125
$ws.Cells.Item(23, 6).Value = 1.079250344657621
$ws.Cells.Item(23, 9).Value = 1.054250312547442
$ws.Cells.Item(23, 10).Value = 1.07101436796918
$ws.Cells.Item(23, 11).Value = 1.070373976223813
$ws.Cells.Item(23, 12).Value = 1.072999269007489
$ws.Cells.Item(23, 13).Value = 1.082734544416927
$ws.Cells.Item(23, 14).Value = 1.072535331218761

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.066188231782312
$ws.Cells.Item(24, 4).Value = 1.06814680567832
$ws.Cells.Item(24, 5).Value = 1.070934805727396
$ws.Cells.Item(24, 6).Value = 1.080746776585388
$ws.Cells.Item(24, 9).Value = 1.054761815858706
$ws.Cells.Item(24, 10).Value = 1.072181426638823
$ws.Cells.Item(24, 11).Value = 1.071408390232673
$ws.Cells.Item(24, 12).Value = 1.074187318835403
$ws.Cells.Item(24, 13).Value = 1.083967771367468
$ws.Cells.Item(24, 14).Value = 1.073704047245578

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.068129327480412
$ws.Cells.Item(25, 4).Value = 1.069657043806171
$ws.Cells.Item(25, 5).Value = 1.072624450310978
$ws.Cells.Item(25, 6).Value = 1.082485998507096
$ws.Cells.Item(25, 9).Value = 1.055350665204271
$ws.Cells.Item(25, 10).Value = 1.073534252787949
$ws.Cells.Item(25, 11).Value = 1.072606838180301
$ws.Cells.Item(25, 12).Value = 1.075565576842517
$ws.Cells.Item(25, 13).Value = 1.085398690167487
$ws.Cells.Item(25, 14).Value = 1.075058794563008

